$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (rows 2 and 3, column B, on both the "zh-cn" and "de-de" sheets)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B2").Value = $newStatus
    $ws.Range("B3").Value = $newStatus
}

# Hyperlink-text colour used elsewhere in the workbook (BGR packed int for
# RGB #6495ED, matching the existing "HyperLink" font/style).
$linkColor = 15570276

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: add "Latest Target File" (E) / "Latest Handback File" (F)
#    entries for rows 2 and 3, plus the handback timestamp in column G.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 (3120aef5-5742-44c4-bfc8-f48d3381e7be)
$wsZh.Range("E2").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/77604d18a2a9f6a23abaf8d8aa5e8165c848e85b/e2e/3120aef5-5742-44c4-bfc8-f48d3381e7be.md", "", "", "3120aef5-5742-44c4-bfc8-f48d3381e7be.md")
$wsZh.Range("E2").Font.Underline = $true
$wsZh.Range("E2").Font.Color = $linkColor
$wsZh.Range("E2").Font.Name = "Calibri"

$wsZh.Range("F2").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e82579785b7537abf9bed3a2b1419f30d2cdc77/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.zh-cn.xlf", "", "", "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.zh-cn.xlf")
$wsZh.Range("F2").Font.Underline = $true
$wsZh.Range("F2").Font.Color = $linkColor
$wsZh.Range("F2").Font.Name = "Calibri"

$wsZh.Range("G2").Value = "2016-01-17 03:19:00"

# Row 3 (a1eefbe0-7ce3-406d-afc9-ca13d215af41)
$wsZh.Range("E3").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/77604d18a2a9f6a23abaf8d8aa5e8165c848e85b/e2e/a1eefbe0-7ce3-406d-afc9-ca13d215af41.md", "", "", "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md")
$wsZh.Range("E3").Font.Underline = $true
$wsZh.Range("E3").Font.Color = $linkColor
$wsZh.Range("E3").Font.Name = "Calibri"

$wsZh.Range("F3").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e82579785b7537abf9bed3a2b1419f30d2cdc77/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.zh-cn.xlf", "", "", "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.zh-cn.xlf")
$wsZh.Range("F3").Font.Underline = $true
$wsZh.Range("F3").Font.Color = $linkColor
$wsZh.Range("F3").Font.Name = "Calibri"

$wsZh.Range("G3").Value = "2016-01-17 03:19:00"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same pattern as zh-cn, with de-de specific filenames/urls.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2 (3120aef5-5742-44c4-bfc8-f48d3381e7be)
$wsDe.Range("E2").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/77604d18a2a9f6a23abaf8d8aa5e8165c848e85b/e2e/3120aef5-5742-44c4-bfc8-f48d3381e7be.md", "", "", "3120aef5-5742-44c4-bfc8-f48d3381e7be.md")
$wsDe.Range("E2").Font.Underline = $true
$wsDe.Range("E2").Font.Color = $linkColor
$wsDe.Range("E2").Font.Name = "Calibri"

$wsDe.Range("F2").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac0cf3b806a666cb8f2010afd452fd5c75235c65/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.de-de.xlf", "", "", "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.de-de.xlf")
$wsDe.Range("F2").Font.Underline = $true
$wsDe.Range("F2").Font.Color = $linkColor
$wsDe.Range("F2").Font.Name = "Calibri"

$wsDe.Range("G2").Value = "2016-01-17 03:19:17"

# Row 3 (a1eefbe0-7ce3-406d-afc9-ca13d215af41)
$wsDe.Range("E3").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/77604d18a2a9f6a23abaf8d8aa5e8165c848e85b/e2e/a1eefbe0-7ce3-406d-afc9-ca13d215af41.md", "", "", "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md")
$wsDe.Range("E3").Font.Underline = $true
$wsDe.Range("E3").Font.Color = $linkColor
$wsDe.Range("E3").Font.Name = "Calibri"

$wsDe.Range("F3").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac0cf3b806a666cb8f2010afd452fd5c75235c65/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.de-de.xlf", "", "", "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.de-de.xlf")
$wsDe.Range("F3").Font.Underline = $true
$wsDe.Range("F3").Font.Color = $linkColor
$wsDe.Range("F3").Font.Name = "Calibri"

$wsDe.Range("G3").Value = "2016-01-17 03:19:17"
